$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark.
#    It currently sits (collapsed) just before the "4.4  Imagery Award"
#    heading; it needs to move into the middle of the sentence
#    "...will start the game in their respective...", splitting that
#    run into " will start the" + " game in their respective ".
# ---------------------------------------------------------------------

# Remove the bookmark from its old location.
$old = $d.Bookmarks.Item("_GoBack")
$old.Delete()

# Find the sentence that must be split and compute the split point
# (right after "...will start the", before " game...").
$rng = $d.Content
$rng.Find.Execute(" will start the game in their respective ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitOffset = " will start the".Length
$splitPoint = $rng.Start + $splitOffset

# Re-add the (collapsed) bookmark at the new location - Word will emit
# this as a run split with <w:bookmarkStart/><w:bookmarkEnd/> in between.
$newRng = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $newRng)

# ---------------------------------------------------------------------
# 2) Remove the "FIRST logo" picture (the floating/anchored picture
#    named "Picture 4") from the page footer. The other footer image
#    (the Broncobots logo) is left in place.
# ---------------------------------------------------------------------

$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)

for ($i = $ftr.Shapes.Count; $i -ge 1; $i--) {
    $shp = $ftr.Shapes.Item($i)
    if ($shp.Name -eq "Picture 4") {
        $shp.Delete()
    }
}

Write-Output "done"
